# Row 25 ("AMD Ryzen 9 5950X") is moved to the end of the existing data
# block (row 39), and every row below it (26-39) shifts up by one row to
# fill the gap. Net effect: rows 25-38 now hold what used to be rows
# 26-39, and row 39 holds what used to be row 25.
#
# Row/range Cut+Insert did not reliably shift cells in this COM host, so
# the move is performed by copying cell values directly, working top to
# bottom so each source row is read before it gets overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 25
$lastRow  = 39

# Remember the row that is being displaced so it can be re-appended at
# the bottom once everything above it has shifted up.
$savedA = $ws.Cells.Item($firstRow, 1).Value2
$savedB = $ws.Cells.Item($firstRow, 2).Value2
$savedC = $ws.Cells.Item($firstRow, 3).Value2

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r - 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r - 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r - 1, 3).Value2 = $ws.Cells.Item($r, 3).Value2
}

$ws.Cells.Item($lastRow, 1).Value2 = $savedA
$ws.Cells.Item($lastRow, 2).Value2 = $savedB
$ws.Cells.Item($lastRow, 3).Value2 = $savedC
